# edit.ps1 - apply LM_UserStories.docx changes via Word COM-interop (PowerShell style)
# Strategy: each target paragraph is located via Find, then its whole <w:p> is
# reconstructed (with w:proofErr spell/gram-check markers split across runs) and
# spliced in via Range.InsertXML at the collapsed start of the paragraph range,
# which replaces that paragraph's content/XML in this engine.

$d = $word.ActiveDocument

function Replace-Paragraph($AnchorText, $NewParagraphXml) {
    $r = $d.Content
    $found = $r.Find.Execute($AnchorText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Anchor text not found"
    }
    $oldPara = $r.Paragraphs(1).Range.Duplicate()
    $insPoint = $oldPara.Duplicate()
    $insPoint.Collapse(1)
    $insPoint.InsertXML($NewParagraphXml)
}

# I Want to ... (je veux...) -> split "Want" with proofErr spellStart/spellEnd
$anchor = 'I Want to ... (je veux'
$newXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">I </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Want</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> to ... (je veux…)</w:t></w:r></w:p>
'@
Replace-Paragraph $anchor $newXml

# so that ... (afin de ...) -> split "so" and "that" with proofErr spell/gram markers
$anchor = 'so that ... (afin de'
$newXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>so</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>that</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> ... (afin de …) </w:t></w:r></w:p>
'@
Replace-Paragraph $anchor $newXml

# Priority -> split "Priority" with proofErr spellStart/spellEnd, keep footnote run
$anchor = 'Priority '
$newXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:proofErr w:type="spellStart"/><w:r><w:t>Priority</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rStyle w:val="Appelnotedebasdep"/></w:rPr><w:footnoteReference w:id="2"/></w:r></w:p>
'@
Replace-Paragraph $anchor $newXml

# CA5 paragraph -> strip _GoBack bookmarkStart/bookmarkEnd
$anchor = 'CA5 : Je veux pouvoir utiliser un compte'
$newXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>CA5 : Je veux pouvoir utiliser un compte externe pour créer mon compte dans l’application (Facebook, Google)</w:t></w:r></w:p>
'@
Replace-Paragraph $anchor $newXml

# CA2 paypal paragraph -> split "paypal" with proofErr spellStart/spellEnd
$anchor = 'CA2 : Je veux pouvoir effectuer un paiem'
$newXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">CA2 : Je veux pouvoir effectuer un paiement avec divers moyen (carte de crédit, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>paypal</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p>
'@
Replace-Paragraph $anchor $newXml

# CA 1 Stripe paragraph -> split "Stripe" with proofErr spellStart/spellEnd
$anchor = 'CA 1 : l’API de Stripe fonctionne grâce '
$newXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">CA 1 : l’API de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Stripe</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> fonctionne grâce au test effectué</w:t></w:r></w:p>
'@
Replace-Paragraph $anchor $newXml

# Final empty paragraphs after the table: the 2nd one gains the _GoBack bookmark pair
# (bookmark previously lived at the end of the CA5 text, now removed above).
$lastParaRange = $d.Paragraphs($d.Paragraphs.Count).Range.Duplicate()
$insPoint2 = $lastParaRange.Duplicate()
$insPoint2.Collapse(1)
$finalXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@
$insPoint2.InsertXML($finalXml)

Write-Output "Edits applied."
